$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new trade row (row 6) with the latest scan result, reusing the
# formatting of the prior trade row so date / boolean styles carry over.
$ws.Range("A5:I5").Copy($ws.Range("A6:I6"))

$ws.Range("A6").Value = 42647.680590277778
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = 9864.9500000000007
$ws.Range("D6").Value = 9874.82
$ws.Range("E6").Value = 104.06
$ws.Range("F6").Value = 104.27
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = 0.2
$ws.Range("I6").Value = $false
